$d = $word.ActiveDocument

$d.Content.Find.Execute("37÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "19÷9=", 2) | Out-Null
$d.Content.Find.Execute("56÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "58÷6=", 2) | Out-Null
$d.Content.Find.Execute("12÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "60÷2=", 2) | Out-Null
$d.Content.Find.Execute("89÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "77÷5=", 2) | Out-Null
$d.Content.Find.Execute("99÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "10÷4=", 2) | Out-Null
$d.Content.Find.Execute("81÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "32÷8=", 2) | Out-Null
$d.Content.Find.Execute("26÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "94÷2=", 2) | Out-Null
$d.Content.Find.Execute("56÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "96÷6=", 2) | Out-Null
$d.Content.Find.Execute("27÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "37÷3=", 2) | Out-Null
$d.Content.Find.Execute("96÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "82÷6=", 2) | Out-Null
$d.Content.Find.Execute("69÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "55÷2=", 2) | Out-Null
$d.Content.Find.Execute("89÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "77÷2=", 2) | Out-Null
$d.Content.Find.Execute("12÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "54÷6=", 2) | Out-Null
$d.Content.Find.Execute("44÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "76÷6=", 2) | Out-Null
$d.Content.Find.Execute("17÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "28÷9=", 2) | Out-Null
$d.Content.Find.Execute("74÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "97÷3=", 2) | Out-Null
$d.Content.Find.Execute("59÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "56÷8=", 2) | Out-Null
$d.Content.Find.Execute("93÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "81÷4=", 2) | Out-Null
$d.Content.Find.Execute("17÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "70÷9=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "13÷3=", 2) | Out-Null
$d.Content.Find.Execute("46÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "53÷2=", 2) | Out-Null
$d.Content.Find.Execute("68÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "54÷9=", 2) | Out-Null
$d.Content.Find.Execute("70÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "82÷2=", 2) | Out-Null
$d.Content.Find.Execute("54÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "75÷9=", 2) | Out-Null
$d.Content.Find.Execute("63÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "23÷8=", 2) | Out-Null
